# This workbook contains one sheet of weekly "Espárragos" (asparagus) price
# records for "Vega Monumental Concepción". The commit re-shuffles the
# per-row data (date, color/variedad, volume, prices, unit, origin) across
# the existing data rows (2-11), while columns A, B, C, E, F, G, I, Q, R
# (market id, market name, region, codreg, category id/name, calidad,
# kg-o-unidades, clasificación) stay identical on every row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: target row -> the (D,H,J,K,L,M,N,O,P) values that must end up there.
$rows = @(
    @{ Row = 2;  D = 44477; H = "Sin especificar"; J = 500; K = 1400; L = 1500; M = 1460; N = '$/kilo';    O = "Provincia de Linares"; P = 1460 },
    @{ Row = 3;  D = 44524; H = "Sin especificar"; J = 200; K = 1500; L = 1600; M = 1550; N = '$/kilo';    O = "Provincia de Talca";   P = 1550 },
    @{ Row = 4;  D = 44526; H = "Sin especificar"; J = 100; K = 1500; L = 1600; M = 1550; N = '$/kilo';    O = "Provincia de Linares"; P = 1550 },
    @{ Row = 5;  D = 44468; H = "Verde";           J = 500; K = 1800; L = 2000; M = 1920; N = '$/kilo';    O = "Provincia de Linares"; P = 1920 },
    @{ Row = 6;  D = 44519; H = "Sin especificar"; J = 250; K = 1200; L = 1300; M = 1240; N = '$/kilo';    O = "Provincia de Linares"; P = 1240 },
    @{ Row = 7;  D = 44496; H = "Sin especificar"; J = 550; K = 1500; L = 2000; M = 1773; N = '$/paquete'; O = "Provincia de Linares"; P = 1773 },
    @{ Row = 8;  D = 44511; H = "Sin especificar"; J = 600; K = 1300; L = 1400; M = 1350; N = '$/kilo';    O = "Provincia de Linares"; P = 1350 },
    @{ Row = 9;  D = 44545; H = "Sin especificar"; J = 550; K = 1700; L = 1800; M = 1755; N = '$/kilo';    O = "Provincia de Linares"; P = 1755 },
    @{ Row = 10; D = 44489; H = "Sin especificar"; J = 600; K = 1400; L = 1500; M = 1450; N = '$/kilo';    O = "Provincia de Linares"; P = 1450 },
    @{ Row = 11; D = 44510; H = "Sin especificar"; J = 600; K = 1300; L = 1400; M = 1350; N = '$/kilo';    O = "Provincia de Linares"; P = 1350 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Range("D$row").Value = $r.D
    $ws.Range("H$row").Value = $r.H
    $ws.Range("J$row").Value = $r.J
    $ws.Range("K$row").Value = $r.K
    $ws.Range("L$row").Value = $r.L
    $ws.Range("M$row").Value = $r.M
    $ws.Range("N$row").Value = $r.N
    $ws.Range("O$row").Value = $r.O
    $ws.Range("P$row").Value = $r.P
}
